$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 27 holds the "END" sentinel value in column A only.
# New enemy categories are being inserted before it, so "END" moves to row 33.
$endValue = $ws.Range("A27").Value2

# New rows of data: category header row (NEW_CATEGORY / name), then
# chance row (100 / dialogue line), repeated for each new enemy.
# Column A values are set first for all three enemies (matching how the
# shared string table was originally populated), then column B values.
$ws.Cells.Item(27, 1).Value = "NEW_CATEGORY"
$ws.Cells.Item(28, 1).Value = 100
$ws.Cells.Item(29, 1).Value = "NEW_CATEGORY"
$ws.Cells.Item(30, 1).Value = 100
$ws.Cells.Item(31, 1).Value = "NEW_CATEGORY"
$ws.Cells.Item(32, 1).Value = 100

$ws.Cells.Item(27, 2).Value = "ijiraq_1"
$ws.Cells.Item(29, 2).Value = "changeling_1"
$ws.Cells.Item(31, 2).Value = "doppelganger_1"

$ws.Cells.Item(28, 2).Value = "THE ARCTIC WINDS ARE HOWLING . . ."
$ws.Cells.Item(30, 2).Value = "NananaNaNAH! Tag, you're it!"
$ws.Cells.Item(32, 2).Value = "I AM SASUN, GOD OF DESTRUCTION."

# Re-append the END sentinel row after the newly inserted rows.
$ws.Cells.Item(33, 1).Value = $endValue
